# Migliorata funzione di wordcloud
#
# The underlying association-rule generation (antecedents/consequents for
# the wordcloud) was fixed upstream; for several rule pairs the
# antecedent/consequent roles - and therefore every metric that depends on
# which side is the antecedent (antecedent/consequent support, confidence,
# lift, leverage, conviction) - come out swapped compared to before. This
# patches the exported "rules" sheet to match the corrected values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row {
    param($Sheet, [int]$Row, [hashtable]$Values)
    foreach ($col in $Values.Keys) {
        $Sheet.Range("$col$Row").Value2 = $Values[$col]
    }
}

# Rows 2 / 3 - frozenset({'2021'}) <-> frozenset({'march'}) swap antecedent/consequent
Set-Row $ws 2 @{ B = "frozenset({'march'})";  C = "frozenset({'2021'})";  D = 0.03314735805710535; E = 0.02297341647522153; G = 0.3366336633663366;  J = 1.4728311184478 }
Set-Row $ws 3 @{ B = "frozenset({'2021'})";   C = "frozenset({'march'})"; D = 0.02297341647522153; E = 0.03314735805710535; G = 0.4857142857142857;  J = 1.879991248222295 }

# Rows 10 / 11 - frozenset({'case'}) <-> frozenset({'new'}) swap antecedent/consequent
Set-Row $ws 10 @{ B = "frozenset({'case'})"; C = "frozenset({'new'})";  D = 0.05874630784378077; E = 0.0653101411224155;   G = 0.3519553072625698; H = 5.388984026276634; J = 1.44232314430248 }
Set-Row $ws 11 @{ B = "frozenset({'new'})";  C = "frozenset({'case'})"; D = 0.0653101411224155;   E = 0.05874630784378077; G = 0.3165829145728643; H = 5.388984026276635; J = 1.377275623081527 }

# Rows 39 / 40 - frozenset({'relief'}) <-> frozenset({'senat'}) swap antecedent/consequent
Set-Row $ws 39 @{ B = "frozenset({'senat'})";  C = "frozenset({'relief'})"; D = 0.02067607482769938; E = 0.02888086642599278; G = 0.5555555555555556; J = 2.185018050541517 }
Set-Row $ws 40 @{ B = "frozenset({'relief'})"; C = "frozenset({'senat'})"; D = 0.02888086642599278; E = 0.02067607482769938; G = 0.3977272727272727; J = 1.626047271984197 }

# Rows 45 / 46 - frozenset({'senat'}) <-> frozenset({'trillion'}) swap antecedent/consequent
Set-Row $ws 45 @{ B = "frozenset({'senat'})";   C = "frozenset({'trillion'})"; D = 0.02067607482769938; E = 0.01017394158188382; G = 0.3174603174603174;  H = 31.20327700972862; J = 1.450210271635845 }
Set-Row $ws 46 @{ B = "frozenset({'trillion'})"; C = "frozenset({'senat'})";   D = 0.01017394158188382; E = 0.02067607482769938; G = 0.6451612903225807;  H = 31.20327700972863; J = 2.75991288003103 }

# Rows 49-51 - {'new','death'}/{'new','case'}/{'death','case'} 2-item antecedent
# group: the antecedent/consequent roles rotate (49<->51 swap in full; row 50's
# numbers were already correct and only the set's display order changes).
Set-Row $ws 49 @{ B = "frozenset({'case', 'death'})"; C = "frozenset({'new'})";  D = 0.01444043321299639; E = 0.0653101411224155;  G = 0.4318181818181818; H = 6.61180904522613;   I = 0.005292534883693376; J = 1.645054151624549 }
Set-Row $ws 50 @{ B = "frozenset({'case', 'new'})" }
Set-Row $ws 51 @{ B = "frozenset({'death', 'new'})"; C = "frozenset({'case'})"; D = 0.009517558254020347; E = 0.05874630784378077; G = 0.6551724137931035; H = 11.15257175881333; I = 0.005676520207591191; J = 2.729635707253037 }

# Rows 52-54 - {'new','report'}/{'new','case'}/{'report','case'} group, same pattern
Set-Row $ws 52 @{ B = "frozenset({'case', 'report'})"; C = "frozenset({'new'})";  D = 0.01247128322940597;  E = 0.0653101411224155;  G = 0.5263157894736842; H = 8.058714625760381;  I = 0.005749332010944605; J = 1.973234146519345 }
Set-Row $ws 53 @{ B = "frozenset({'case', 'new'})" }
Set-Row $ws 54 @{ B = "frozenset({'report', 'new'})"; C = "frozenset({'case'})"; D = 0.008861174926156876; E = 0.05874630784378077; G = 0.7407407407407407; H = 12.60914545830747; I = 0.00604327196856512; J = 3.630549955459702 }

# Rows 55-58 - {'trillion','relief'}/{'trillion','senat'}/{'relief','senat'}/{'trillion'} group, same pattern
Set-Row $ws 55 @{ B = "frozenset({'senat', 'relief'})";  C = "frozenset({'trillion'})"; D = 0.01148670823761076;  E = 0.01017394158188382;  G = 0.4571428571428571; H = 44.93271889400922; I = 0.005134201524330182; J = 1.823363791822845 }
Set-Row $ws 56 @{ B = "frozenset({'senat', 'trillion'})" }
Set-Row $ws 57 @{ B = "frozenset({'relief', 'trillion'})"; C = "frozenset({'senat'})"; D = 0.006563833278634723; E = 0.02067607482769938; G = 0.7999999999999999; H = 38.69206349206349; I = 0.005115352314882183; J = 4.896619625861502 }
Set-Row $ws 58 @{ C = "frozenset({'senat', 'relief'})" }
